$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 14 content: Manufacturer / Mfg Part # / Description / Package stays the same
$ws.Range("D14").Value = "Texas Instruments"
$ws.Range("E14").Value = "SN74HCS125QBQARQ1"
$ws.Range("F14").Value = "Automotive Schmitt-trigger inputs quadruple bus buffer gates with 3-state outputs 14-WQFN -40 to 125"

# Adjust column F width (closest achievable value to the target 87.6640625
# given this engine's column-width rounding granularity)
$ws.Columns.Item(6).ColumnWidth = 86.83

# Update view selection to match the edited cell
$ws.Range("F14").Select()
